$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "the day after tomorrow"
$ws.Range("B2").Value = "あさって"
$ws.Range("A3").Value = "rain"
$ws.Range("B3").Value = "雨|あめ"
$ws.Range("A4").Value = "office worker"
$ws.Range("B4").Value = "会社員|かいしゃいん"
$ws.Range("A5").Value = "camera"
$ws.Range("B5").Value = "カメラ"
$ws.Range("A6").Value = "karaoke"
$ws.Range("B6").Value = "カラオケ"
$ws.Range("A7").Value = "air"
$ws.Range("B7").Value = "空気|くうき"
$ws.Range("A8").Value = "this morning"
$ws.Range("B8").Value = "今朝|けさ"
$ws.Range("A9").Value = "blackboard"
$ws.Range("B9").Value = "黒板|こくばん"
$ws.Range("A10").Value = "this month"
$ws.Range("B10").Value = "今月|こんげつ"
$ws.Range("A11").Value = "job; work; occupation"
$ws.Range("B11").Value = "仕事|しごと"
$ws.Range("A12").Value = "college student"
$ws.Range("B12").Value = "大学生|だいがくせい"
$ws.Range("A13").Value = "weather forecast"
$ws.Range("B13").Value = "天気予報|てんきよほう"
$ws.Range("A14").Value = "place"
$ws.Range("B14").Value = "所|ところ"
$ws.Range("A15").Value = "tomato"
$ws.Range("B15").Value = "トマト"
$ws.Range("A16").Value = "summer"
$ws.Range("B16").Value = "夏|なつ"
$ws.Range("A17").Value = "something"
$ws.Range("B17").Value = "何か|なにか"
$ws.Range("A18").Value = "party"
$ws.Range("B18").Value = "パーティー"
$ws.Range("A19").Value = "barbecue"
$ws.Range("B19").Value = "バーベキュー"
$ws.Range("A20").Value = "chopsticks"
$ws.Range("B20").Value = "はし"
$ws.Range("A21").Value = "winter"
$ws.Range("B21").Value = "冬|ふゆ"
$ws.Range("A22").Value = "homestay; living with a local family"
$ws.Range("B22").Value = "ホームステイ"
$ws.Range("A23").Value = "every week"
$ws.Range("B23").Value = "毎週|まいしゅう"
$ws.Range("A24").Value = "next month"
$ws.Range("B24").Value = "来月|らいげつ"
$ws.Range("A25").Value = "skillful; good at"
$ws.Range("B25").Value = "上手|じょうず(な)"
$ws.Range("A26").Value = "clumsy; poor at"
$ws.Range("B26").Value = "下手|へた(な)"
$ws.Range("A27").Value = "famous"
$ws.Range("B27").Value = "有名|ゆうめい(な)"
$ws.Range("A28").Value = "it rains"
$ws.Range("B28").Value = "雨が降る|あめがふる"
$ws.Range("A29").Value = "to wash"
$ws.Range("B29").Value = "洗う|あらう"
$ws.Range("A30").Value = "to say"
$ws.Range("B30").Value = "言う|いう"
$ws.Range("A31").Value = "to need"
$ws.Range("B31").Value = "いる"
$ws.Range("A32").Value = "to be late"
$ws.Range("B32").Value = "遅くなる|おそくなる"
$ws.Range("A33").Value = "to think"
$ws.Range("B33").Value = "思う|おもう"
$ws.Range("A34").Value = "to cut"
$ws.Range("B34").Value = "切る|きる"
$ws.Range("A35").Value = "to make"
$ws.Range("B35").Value = "作る|つくる"
$ws.Range("A36").Value = "to take (a thing)"
$ws.Range("B36").Value = "持っていく|もっていく"
$ws.Range("A37").Value = "to stare (at...)"
$ws.Range("B37").Value = "じろじろ見る|じろじろみる"
$ws.Range("A38").Value = "to throw away"
$ws.Range("B38").Value = "捨てる|すてる"
$ws.Range("A39").Value = "to begin"
$ws.Range("B39").Value = "始める|はじめる"
$ws.Range("A40").Value = "to drive"
$ws.Range("B40").Value = "運転する|うんてんする"
$ws.Range("A41").Value = "to do laundry"
$ws.Range("B41").Value = "洗濯する|せんたくする"
$ws.Range("A42").Value = "to clean"
$ws.Range("B42").Value = "掃除する|そうじする"
$ws.Range("A43").Value = "to call"
$ws.Range("B43").Value = "電話する|でんわする"
$ws.Range("A44").Value = "to cook"
$ws.Range("B44").Value = "料理する|りょうりする"
$ws.Range("A45").Value = "always"
$ws.Range("B45").Value = "いつも"
$ws.Range("A46").Value = "uh-uh; no"
$ws.Range("B46").Value = "ううん"
$ws.Range("A47").Value = "uh-huh; yes"
$ws.Range("B47").Value = "うん"
$ws.Range("A48").Value = "Cheers! (a toast)"
$ws.Range("B48").Value = "乾杯|かんぱい"
$ws.Range("A49").Value = "That's too bad."
$ws.Range("B49").Value = "残念(ですね)|ざんねん(ですね)"
$ws.Range("A50").Value = "about...; concerning..."
$ws.Range("B50").Value = "～について"
$ws.Range("A51").Value = "not...yet"
$ws.Range("B51").Value = "まだ+negative"
$ws.Range("A52").Value = "all (of the people) together"
$ws.Range("B52").Value = "みんなで"
$ws.Range("A53").Value = "Rice"
$ws.Range("B53").Value = "ご飯|ごはん"
$ws.Range("A54").Value = "Miso soup"
$ws.Range("B54").Value = "みそ汁|みそしる"
$ws.Range("A55").Value = "side dish"
$ws.Range("B55").Value = "おかず"
$ws.Range("A56").Value = "set meal"
$ws.Range("B56").Value = "定食|ていしょく"
$ws.Range("A57").Value = "Curry with rice"
$ws.Range("B57").Value = "カレーライス"
$ws.Range("A58").Value = "Deep-fried shrimp"
$ws.Range("B58").Value = "えびフライ"
$ws.Range("A59").Value = "Ramen noodles"
$ws.Range("B59").Value = "ラーメン"
$ws.Range("A60").Value = "Udon noodles"
$ws.Range("B60").Value = "うどん"
$ws.Range("A61").Value = "Spaghetti"
$ws.Range("B61").Value = "スパゲッティ"
$ws.Range("A62").Value = "Dumplings"
$ws.Range("B62").Value = "ぎょうざ"
$ws.Range("A63").Value = "Beef rice bowl"
$ws.Range("B63").Value = "牛丼|ぎゅうどん"
$ws.Range("A64").Value = "Hamburger steak"
$ws.Range("B64").Value = "ハンバーガ"
$ws.Range("A65").Value = "Raw seafood"
$ws.Range("B65").Value = "さしみ"
$ws.Range("A66").Value = "Savory pancake"
$ws.Range("B66").Value = "お好み焼き|おこのみやき"
$ws.Range("A67").Value = "Toast"
$ws.Range("B67").Value = "トースト"
$ws.Range("A68").Value = "Soup"
$ws.Range("B68").Value = "スープ"
$ws.Range("A69").Value = "Yogurt"
$ws.Range("B69").Value = "ヨーグルト"
$ws.Range("A70").Value = "Broiled fish"
$ws.Range("B70").Value = "焼き魚|やきざかな"
$ws.Range("A71").Value = "Egg"
$ws.Range("B71").Value = "たまご"
$ws.Range("A72").Value = "office worker"
$ws.Range("B72").Value = "会社員|かいしゃいん"
$ws.Range("A73").Value = "store clerk"
$ws.Range("B73").Value = "店員|てんいん"
$ws.Range("A74").Value = "member"
$ws.Range("B74").Value = "会員|かいいん"
$ws.Range("A75").Value = "station staff"
$ws.Range("B75").Value = "駅員|えきいん"
$ws.Range("A76").Value = "new"
$ws.Range("B76").Value = "新しい|あたらしい"
$ws.Range("A77").Value = "newspaper"
$ws.Range("B77").Value = "新聞|しんぶん"
$ws.Range("A78").Value = "Bullet Train"
$ws.Range("B78").Value = "新幹線|しんかんせん"
$ws.Range("A79").Value = "fresh"
$ws.Range("B79").Value = "新鮮な|しんせんな"
$ws.Range("A80").Value = "to listen"
$ws.Range("B80").Value = "聞く|きく"
$ws.Range("A81").Value = "can be heard"
$ws.Range("B81").Value = "聞こえる|きこえる"
$ws.Range("A82").Value = "to make"
$ws.Range("B82").Value = "作る|つくる"
$ws.Range("A83").Value = "composition"
$ws.Range("B83").Value = "作文|さくぶん"
$ws.Range("A84").Value = "artistic piece"
$ws.Range("B84").Value = "作品|さくひん"
$ws.Range("A85").Value = "author"
$ws.Range("B85").Value = "作者|さくしゃ"
$ws.Range("A86").Value = "job"
$ws.Range("B86").Value = "仕事|しごと"
$ws.Range("A87").Value = "revenge"
$ws.Range("B87").Value = "仕返し|しかえし"
$ws.Range("A88").Value = "to serve; to work under"
$ws.Range("B88").Value = "仕える|つかえる"
$ws.Range("A89").Value = "thing"
$ws.Range("B89").Value = "事|こと"
$ws.Range("A90").Value = "fire"
$ws.Range("B90").Value = "火事|かじ"
$ws.Range("A91").Value = "meal"
$ws.Range("B91").Value = "食事|しょくじ"
$ws.Range("A92").Value = "reply"
$ws.Range("B92").Value = "返事|へんじ"
$ws.Range("A93").Value = "train"
$ws.Range("B93").Value = "電車|でんしゃ"
$ws.Range("A94").Value = "electricity"
$ws.Range("B94").Value = "電気|でんき"
$ws.Range("A95").Value = "telephone"
$ws.Range("B95").Value = "電話|でんわ"
$ws.Range("A96").Value = "battery"
$ws.Range("B96").Value = "電池|でんち"
$ws.Range("A97").Value = "electronic dictionary"
$ws.Range("B97").Value = "電子辞書|でんしじしょ"
$ws.Range("A98").Value = "car"
$ws.Range("B98").Value = "車|くるま"
$ws.Range("A99").Value = "bicycle"
$ws.Range("B99").Value = "自転車|じてんしゃ"
$ws.Range("A100").Value = "wheel chair"
$ws.Range("B100").Value = "車いす|くるまいす"
$ws.Range("A101").Value = "parking lot"
$ws.Range("B101").Value = "駐車場|ちゅうしゃじょう"
$ws.Range("A112").Value = "to think"
$ws.Range("B112").Value = "思う|おもう"
$ws.Range("A113").Value = "mysterious"
$ws.Range("B113").Value = "不思議な|ふしぎな"
$ws.Range("A114").Value = "to recall; to remember"
$ws.Range("B114").Value = "思い出す|おもいだす"
$ws.Range("A115").Value = "next"
$ws.Range("B115").Value = "次|つぎ"
$ws.Range("A116").Value = "second daughter"
$ws.Range("B116").Value = "次女|じじょ"
$ws.Range("A117").Value = "table of contents"
$ws.Range("B117").Value = "目次|もくじ"
$ws.Range("A118").Value = "next time"
$ws.Range("B118").Value = "次回|じかい"
$ws.Range("A119").Value = "what"
$ws.Range("B119").Value = "何|なに"
$ws.Range("A120").Value = "what time"
$ws.Range("B120").Value = "何時|なんじ"
$ws.Range("A122").Value = "something"
$ws.Range("B122").Value = "何か|なにか"
